# Apply the "Uploading newest copy of EPS US-develop" edits to the
# Transportation Technology Logit Exponents workbook.
#
# Sheets: "About" (1), "A54.tranSubsector_logit_revised" (2), "TTLE" (3)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "About" sheet: update/add explanatory notes
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# B3 stays "Calibration" (unchanged text, no-op but kept for clarity)
$about.Range("B3").Value = "Calibration"

# A13 is a brand-new note (previously empty) - written before A12 so the
# shared-string table append order matches the source file's ordering.
$about.Range("A13").Value = "For non-road we use -0.1 because of costs preventing the model from solving."

# A12 note text changed
$about.Range("A12").Value = "We use calibrated values in onroad sectors."

# A15 wording changed from "Modified Logit" to "Unmodified Logit"
$about.Range("A15").Value = 'For more on this, see the "Unmodified Logit" equation description at:'

# A16 (URL) unchanged
$about.Range("A16").Value = "https://jgcri.github.io/gcam-doc/choice.html"

# ---------------------------------------------------------------------
# 2) "TTLE" sheet: replace the live cross-sheet formulas with hard-coded
#    literal numbers (the linkage to A54.tranSubsector_logit_revised was
#    broken out into plain values, and several values were also revised).
# ---------------------------------------------------------------------
$ttle = $wb.Worksheets.Item("TTLE")

$ttle.Range("B2").Value = -80
$ttle.Range("C2").Value = -40

$ttle.Range("B3").Value = -30
$ttle.Range("C3").Value = -30

$ttle.Range("B4").Value = -0.1
$ttle.Range("C4").Value = -0.1

$ttle.Range("B5").Value = -0.15
$ttle.Range("C5").Value = -0.1

$ttle.Range("B6").Value = -0.1
$ttle.Range("C6").Value = -0.1

$ttle.Range("B7").Value = -40
$ttle.Range("C7").Value = -40

# ---------------------------------------------------------------------
# 3) Selections / active tab: the saved view moves from "About"
#    (cell A12 selected) to "TTLE" being the active/selected tab
#    (cell G3 selected); "About"'s remembered selection becomes A16.
# ---------------------------------------------------------------------
$about.Activate()
$about.Range("A16").Select() | Out-Null

$ttle.Activate()
$ttle.Range("G3").Select() | Out-Null
